$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1582.9259
$ws.Range("I17").Value = 1178.2222
$ws.Range("J17").Value = 1785.2778
$ws.Range("K17").Value = 3534.6666
$ws.Range("L17").Value = 5355.8334
$ws.Range("M17").Value = -3366.6666
$ws.Range("N17").Value = -5691.8334
$ws.Range("H53").Value = 1967.2106
$ws.Range("I53").Value = 285.66666
$ws.Range("K53").Value = 285.66666
$ws.Range("M53").Value = 351.33334
$ws.Range("H62").Value = 7494.5
$ws.Range("I62").Value = 6189.6
$ws.Range("K62").Value = 6189.6
$ws.Range("M62").Value = -5565.6
$ws.Range("H65").Value = 7494.5
$ws.Range("I65").Value = 6189.6
$ws.Range("K65").Value = 30948
$ws.Range("M65").Value = -27828
$ws.Range("H112").Value = 1939.8572
$ws.Range("I112").Value = 1300
$ws.Range("J112").Value = 1971.85
$ws.Range("K112").Value = 3900
$ws.Range("L112").Value = 5915.549999999999
$ws.Range("M112").Value = -2792
$ws.Range("N112").Value = -8131.549999999999
$ws.Range("H116").Value = 8742.615
$ws.Range("I116").Value = 7639.6
$ws.Range("J116").Value = 9432
$ws.Range("K116").Value = 7639.6
$ws.Range("L116").Value = 9432
$ws.Range("M116").Value = -4197.6
$ws.Range("N116").Value = -16316
$ws.Range("H125").Value = 2144.3845
$ws.Range("J125").Value = 2139.8572
$ws.Range("L125").Value = 19258.7148
$ws.Range("N125").Value = -24178.7148
$ws.Range("H127").Value = 2600.7
$ws.Range("I127").Value = 1172.75
$ws.Range("J127").Value = 3552.6667
$ws.Range("K127").Value = 3518.25
$ws.Range("L127").Value = 10658.0001
$ws.Range("M127").Value = 1441.75
$ws.Range("N127").Value = -20578.0001
$ws.Range("H129").Value = 2408.6155
$ws.Range("J129").Value = 2682.9092
$ws.Range("L129").Value = 8048.7276
$ws.Range("N129").Value = -18048.7276
$ws.Range("H138").Value = 5300.384
$ws.Range("J138").Value = 5562.0464
$ws.Range("L138").Value = 16686.1392
$ws.Range("N138").Value = -26966.1392

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9248.98
$ws.Range("I32").Value = 8453.68
$ws.Range("K32").Value = 8453.68
$ws.Range("M32").Value = -8166.68
$ws.Range("H61").Value = 6603.44
$ws.Range("I61").Value = 6855.5
$ws.Range("J61").Value = 6282.636
$ws.Range("K61").Value = 6855.5
$ws.Range("L61").Value = 6282.636
$ws.Range("M61").Value = -6643.5
$ws.Range("N61").Value = -6706.636
$ws.Range("H74").Value = 2640.0688
$ws.Range("I74").Value = 2120.8823
$ws.Range("K74").Value = 2120.8823
$ws.Range("M74").Value = -1246.8823
$ws.Range("H77").Value = 2640.0688
$ws.Range("I77").Value = 2120.8823
$ws.Range("K77").Value = 10604.4115
$ws.Range("M77").Value = -6236.411500000002
$ws.Range("H97").Value = 789.8946999999999
$ws.Range("I97").Value = 593.0833
$ws.Range("K97").Value = 593.0833
$ws.Range("M97").Value = -97.08330000000001
$ws.Range("H132").Value = 2481.6167
$ws.Range("I132").Value = 2070.6545
$ws.Range("K132").Value = 6211.9635
$ws.Range("M132").Value = -3681.9635
$ws.Range("H135").Value = 57342.734
$ws.Range("J135").Value = 57342.734
$ws.Range("L135").Value = 57342.734
$ws.Range("N135").Value = -67482.734
$ws.Range("H136").Value = 6603.44
$ws.Range("I136").Value = 6855.5
$ws.Range("J136").Value = 6282.636
$ws.Range("K136").Value = 20566.5
$ws.Range("L136").Value = 18847.908
$ws.Range("M136").Value = -18016.5
$ws.Range("N136").Value = -23947.908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3099.5
$ws.Range("J20").Value = 3299.5
$ws.Range("L20").Value = 3299.5
$ws.Range("N20").Value = -3793.5
$ws.Range("H94").Value = 2719.6667
$ws.Range("I94").Value = 2705.65
$ws.Range("K94").Value = 2705.65
$ws.Range("M94").Value = -2254.65
$ws.Range("H134").Value = 3919.913
$ws.Range("I134").Value = 3507.2
$ws.Range("J134").Value = 6671.3335
$ws.Range("K134").Value = 10521.6
$ws.Range("L134").Value = 20014.0005
$ws.Range("M134").Value = -7986.599999999999
$ws.Range("N134").Value = -25084.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37435.805
$ws.Range("I31").Value = 2450.3333
$ws.Range("J31").Value = 85877.234
$ws.Range("K31").Value = 2450.3333
$ws.Range("L31").Value = 85877.234
$ws.Range("M31").Value = -2155.3333
$ws.Range("N31").Value = -86467.234
$ws.Range("H34").Value = 37435.805
$ws.Range("I34").Value = 2450.3333
$ws.Range("J34").Value = 85877.234
$ws.Range("K34").Value = 2450.3333
$ws.Range("L34").Value = 85877.234
$ws.Range("M34").Value = -2248.3333
$ws.Range("N34").Value = -86281.234
$ws.Range("H58").Value = 4563.7144
$ws.Range("I58").Value = 1941.0714
$ws.Range("J58").Value = 7186.357
$ws.Range("K58").Value = 1941.0714
$ws.Range("L58").Value = 7186.357
$ws.Range("M58").Value = -1738.0714
$ws.Range("N58").Value = -7592.357
$ws.Range("H132").Value = 4076.14
$ws.Range("I132").Value = 3615.7144
$ws.Range("J132").Value = 6493.375
$ws.Range("K132").Value = 10847.1432
$ws.Range("L132").Value = 19480.125
$ws.Range("M132").Value = -8317.143199999999
$ws.Range("N132").Value = -24540.125
$ws.Range("H134").Value = 3252.5
$ws.Range("I134").Value = 2466.1428
$ws.Range("J134").Value = 8757
$ws.Range("K134").Value = 7398.428400000001
$ws.Range("L134").Value = 26271
$ws.Range("M134").Value = -4863.428400000001
$ws.Range("N134").Value = -31341
$ws.Range("H136").Value = 4563.7144
$ws.Range("I136").Value = 1941.0714
$ws.Range("J136").Value = 7186.357
$ws.Range("K136").Value = 5823.2142
$ws.Range("L136").Value = 21559.071
$ws.Range("M136").Value = -3273.2142
$ws.Range("N136").Value = -26659.071
$ws.Range("H141").Value = 254895.84
$ws.Range("J141").Value = 271137.16
$ws.Range("L141").Value = 271137.16
$ws.Range("N141").Value = -281497.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2084245.9
$ws.Range("J107").Value = 6251600
$ws.Range("L107").Value = 18754800
$ws.Range("N107").Value = -18758640
$ws.Range("H113").Value = 2781.16
$ws.Range("J113").Value = 2916.7273
$ws.Range("L113").Value = 8750.1819
$ws.Range("N113").Value = -13090.1819
$ws.Range("H124").Value = 7897.7144
$ws.Range("I124").Value = 7499
$ws.Range("J124").Value = 7964.1665
$ws.Range("K124").Value = 22497
$ws.Range("L124").Value = 23892.4995
$ws.Range("M124").Value = -17587
$ws.Range("N124").Value = -33712.49950000001
$ws.Range("H136").Value = 3048.4827
$ws.Range("I136").Value = 2956.28
$ws.Range("J136").Value = 3624.75
$ws.Range("K136").Value = 8868.84
$ws.Range("L136").Value = 10874.25
$ws.Range("M136").Value = -3768.84
$ws.Range("N136").Value = -21074.25
$ws.Range("H137").Value = 95203.73
$ws.Range("I137").Value = 2887.5
$ws.Range("J137").Value = 115718.445
$ws.Range("K137").Value = 8662.5
$ws.Range("L137").Value = 347155.335
$ws.Range("M137").Value = -3562.5
$ws.Range("N137").Value = -357355.335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 45000
$ws.Range("I62").Value = 45000
$ws.Range("K62").Value = 45000
$ws.Range("M62").Value = -44314
$ws.Range("H65").Value = 45000
$ws.Range("I65").Value = 45000
$ws.Range("K65").Value = 135000
$ws.Range("M65").Value = -131568
$ws.Range("H126").Value = 4084.1035
$ws.Range("I126").Value = 2966.647
$ws.Range("J126").Value = 5667.1665
$ws.Range("K126").Value = 8899.940999999999
$ws.Range("L126").Value = 17001.4995
$ws.Range("M126").Value = -6429.940999999999
$ws.Range("N126").Value = -21941.4995
$ws.Range("H132").Value = 4374.75
$ws.Range("I132").Value = 2646.9
$ws.Range("K132").Value = 7940.700000000001
$ws.Range("M132").Value = -5410.700000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9886
$ws.Range("I40").Value = 10273.667
$ws.Range("K40").Value = 10273.667
$ws.Range("M40").Value = -10137.667
$ws.Range("H68").Value = 5108.205
$ws.Range("I68").Value = 3321.75
$ws.Range("K68").Value = 3321.75
$ws.Range("M68").Value = -2572.75
$ws.Range("H71").Value = 5108.205
$ws.Range("I71").Value = 3321.75
$ws.Range("K71").Value = 16608.75
$ws.Range("M71").Value = -12864.75
$ws.Range("H131").Value = 129999.5
$ws.Range("J131").Value = 129999.5
$ws.Range("L131").Value = 129999.5
$ws.Range("N131").Value = -140079.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2924.72
$ws.Range("I122").Value = 1995.4
$ws.Range("K122").Value = 5986.200000000001
$ws.Range("M122").Value = -3536.200000000001
$ws.Range("H126").Value = 4068.7693
$ws.Range("I126").Value = 2175.8823
$ws.Range("K126").Value = 6527.646900000001
$ws.Range("M126").Value = -4057.646900000001
$ws.Range("H132").Value = 2548.682
$ws.Range("I132").Value = 1407.6177
$ws.Range("K132").Value = 4222.8531
$ws.Range("M132").Value = -1692.8531
